# Add "The Logic of Political Survival" as a new row in the dataset list.
# The sheet is kept sorted alphabetically by column A (name), and this
# title sorts between row 218 ("The Dyadic Cyber Incident and Dispute
# Data") and the old row 219 ("The Political Constraint Index"), so the
# new entry is inserted as row 219, pushing everything below it down by
# one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 219 (shifts old rows 219:245 down to 220:246),
# then wipe any residual formatting/content the insert carried over from
# the row above so only the cells we fill in below are populated.
$ws.Rows("219:219").Insert()
$ws.Rows("219:219").Clear()

# --- Fill in the new row's data -------------------------------------------------
# (values are assigned in the same order the new shared strings were first
# introduced: name, link, zip file, then topics)
$ws.Range("A219").Value = "The Logic of Political Survival"
$ws.Range("B219").Value = "international relations"
$ws.Range("C219").Value = "http://www.nyu.edu/gsas/dept/politics/data/bdm2s2/Logic.htm"
$ws.Range("F219").Value = 1
$ws.Range("G219").Value = 1
$ws.Range("H219").Value = 1
$ws.Range("I219").Value = 1
$ws.Range("J219").Value = 1
$ws.Range("K219").Value = 1763
$ws.Range("L219").Value = 2010
$ws.Range("M219").Value = "online"
$ws.Range("N219").Value = "no"
$ws.Range("O219").Value = 1
$ws.Range("V219").Value = "http://www.nyu.edu/gsas/dept/politics/data/bdm2s2/bdm2s2_nation_year_data_may2002_webversion.zip"
$ws.Range("W219").Value = "country"
$ws.Range("X219").Value = "year"
$ws.Range("AB219").Value = 20180526
$ws.Range("D219").Value = "leaders, political survival"

# --- Hyperlinks for the link + zip file columns ----------------------------------
$ws.Hyperlinks.Add($ws.Range("C219"), "http://www.nyu.edu/gsas/dept/politics/data/bdm2s2/Logic.htm")
$ws.Hyperlinks.Add($ws.Range("V219"), "http://www.nyu.edu/gsas/dept/politics/data/bdm2s2/bdm2s2_nation_year_data_may2002_webversion.zip")

# Hyperlinks.Add re-styles the cell with a freshly-minted style; put the
# original shared "Hyperlink" cell style back so styling matches the other
# link cells in the sheet instead of growing a duplicate style entry.
$ws.Range("C219").Style = "Hyperlink"
$ws.Range("V219").Style = "Hyperlink"

# --- View state (matches where the author was looking after the edit) -----------
$ws.Range("A219").Select()
$excel.ActiveWindow.ScrollRow = 221
